$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row text corrections (accents removed in a few header labels)
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Correo electronico (*)"
$ws.Range("F1").Value = "Subarea 1 (*)"
$ws.Range("G1").Value = "Subarea 2"
$ws.Range("H1").Value = "Subarea 3"

# ---------------------------------------------------------------------------
# 2. Preserve D2's current cell format (the custom "hyperlink-like" blue
#    font, style index 2) on a scratch cell so it can be re-applied later,
#    since Hyperlinks.Add() below always forces its own builtin Hyperlink
#    style onto the target cell.
# ---------------------------------------------------------------------------
$ws.Range("D2").Copy() | Out-Null
$ws.Range("Q3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Row 2 becomes "Maria Urbina" (was "Pedro Paredes")
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Maria"
$ws.Range("B2").Value = "Urbina"
$ws.Range("C2").Value = "maurbina"
$ws.Range("D2").Value = "murbina@yopmail.com"
$ws.Range("I2").Value = "Femenino"
# E2/F2/H2/K2/L2/J2/N2 stay as-is (Biociencias / Parasitología / Bioquímica /
# "Licenciado " / "Investigación celular" / 45869632 / 4168546321)

# ---------------------------------------------------------------------------
# 4. Row 3 becomes "Emilia Urbina" (was "Mary James")
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Emilia"
$ws.Range("B3").Value = "Urbina"
$ws.Range("C3").Value = "eurbina2"
$ws.Range("D3").Value = "eurbina@yopmail.com"
$ws.Range("I3").Value = "Femenino"
$ws.Range("K3").Value = "Licenciado"
$ws.Range("L3").Value = "Genética aplicada"
# E3/F3/H3 stay as-is (Biociencias / Parasitología / Bioquímica); J3/N3 unchanged

# ---------------------------------------------------------------------------
# 5. New row 4: "Joan Magallanes"
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Joan"
$ws.Range("B4").Value = "Magallanes"
$ws.Range("C4").Value = "jmagallanes"
$ws.Range("D4").Value = "jmagallanes@yopmail.com"
$ws.Range("E4").Value = "Tecnología"
$ws.Range("F4").Value = "Computación"
$ws.Range("I4").Value = "Masculino"
$ws.Range("J4").Value = 156423651
$ws.Range("K4").Value = "Ingeniero"
$ws.Range("L4").Value = "Desarrollo de aplicaciones"
$ws.Range("N4").Value = 4125658574

# ---------------------------------------------------------------------------
# 6. Hyperlinks: drop the old pair (paredesp@/mjames@) and add the new pair
#    on D2 (murbina@yopmail.com) and D4 (jmagallanes@yopmail.com). D3 no
#    longer carries a functional hyperlink (matches target workbook).
#    Hyperlinks.Delete() on a range clears the sheet's hyperlink collection,
#    so do it once and re-add exactly what's needed.
# ---------------------------------------------------------------------------
$ws.Range("D2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:murbina@yopmail.com", "", "", "murbina@yopmail") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:jmagallanes@yopmail.com", "", "", "jmagallanes@yopmail") | Out-Null

# Re-apply the original (non-builtin) blue-text cell format that Hyperlinks.Add
# just overwrote, using the format stashed in the scratch cell, then clean up.
$ws.Range("Q3").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("Q3").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 7. Sheet view bookkeeping (dimension grows to row 4; selection/topLeft reset)
# ---------------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D9").Select() | Out-Null
